$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5163.5
$ws.Range("I86").Value = 4995.75
$ws.Range("K86").Value = 4995.75
$ws.Range("M86").Value = -3872.75
$ws.Range("H89").Value = 5163.5
$ws.Range("I89").Value = 4995.75
$ws.Range("K89").Value = 24978.75
$ws.Range("M89").Value = -19362.75
$ws.Range("H110").Value = 40244.5
$ws.Range("J110").Value = 40244.5
$ws.Range("L110").Value = 40244.5
$ws.Range("N110").Value = -48424.5
$ws.Range("H125").Value = 14227.477
$ws.Range("I125").Value = 4367.846
$ws.Range("K125").Value = 39310.61399999999
$ws.Range("M125").Value = -36850.61399999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 145.6
$ws.Range("I5").Value = 119.76923
$ws.Range("K5").Value = 119.76923
$ws.Range("M5").Value = -7.769229999999993
$ws.Range("H24").Value = 36160
$ws.Range("J24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("N24").Value = -35748
$ws.Range("H34").Value = 284500
$ws.Range("I34").Value = 284500
$ws.Range("K34").Value = 284500
$ws.Range("M34").Value = -284229
$ws.Range("H36").Value = 1289.6
$ws.Range("I36").Value = 1299.5
$ws.Range("J36").Value = 1250
$ws.Range("K36").Value = 1299.5
$ws.Range("L36").Value = 1250
$ws.Range("M36").Value = -953.5
$ws.Range("N36").Value = -1942
$ws.Range("H61").Value = 10006.9
$ws.Range("I61").Value = 10008.75
$ws.Range("K61").Value = 10008.75
$ws.Range("M61").Value = -9796.75
$ws.Range("H100").Value = 36160
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164
$ws.Range("H132").Value = 7619.2383
$ws.Range("I132").Value = 7565.8335
$ws.Range("K132").Value = 22697.5005
$ws.Range("M132").Value = -20167.5005
$ws.Range("H136").Value = 10006.9
$ws.Range("I136").Value = 10008.75
$ws.Range("K136").Value = 30026.25
$ws.Range("M136").Value = -27476.25

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 145.6
$ws.Range("I4").Value = 119.76923
$ws.Range("K4").Value = 119.76923
$ws.Range("M4").Value = -4.769229999999993
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1991
$ws.Range("I132").Value = 2058.4
$ws.Range("K132").Value = 6175.200000000001
$ws.Range("M132").Value = -3645.200000000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 269.28
$ws.Range("I2").Value = 105
$ws.Range("J2").Value = 398.35715
$ws.Range("K2").Value = 630
$ws.Range("L2").Value = 2390.1429
$ws.Range("M2").Value = -517
$ws.Range("N2").Value = -2616.1429
$ws.Range("H11").Value = 1180808.5
$ws.Range("J11").Value = 1667263.8
$ws.Range("L11").Value = 5001791.4
$ws.Range("N11").Value = -5002071.4
$ws.Range("H26").Value = 863.63635
$ws.Range("I26").Value = 1699.8
$ws.Range("K26").Value = 5099.4
$ws.Range("M26").Value = -4811.4
$ws.Range("H33").Value = 499.72726
$ws.Range("J33").Value = 499.77777
$ws.Range("L33").Value = 2998.66662
$ws.Range("N33").Value = -3564.66662
$ws.Range("H74").Value = 60756
$ws.Range("J74").Value = 66500
$ws.Range("L74").Value = 199500
$ws.Range("N74").Value = -201622
$ws.Range("H77").Value = 60756
$ws.Range("J77").Value = 66500
$ws.Range("L77").Value = 598500
$ws.Range("N77").Value = -609108
$ws.Range("H107").Value = 855
$ws.Range("J107").Value = 403.33334
$ws.Range("L107").Value = 1210.00002
$ws.Range("N107").Value = -5050.000019999999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10001000
$ws.Range("I10").Value = 20000000
$ws.Range("K10").Value = 20000000
$ws.Range("M10").Value = -19999831
$ws.Range("H15").Value = 495011800
$ws.Range("J15").Value = 495011800
$ws.Range("L15").Value = 495011800
$ws.Range("N15").Value = -495012376
$ws.Range("H80").Value = 3324.2144
$ws.Range("I80").Value = 2108.5833
$ws.Range("J80").Value = 4235.9375
$ws.Range("K80").Value = 2108.5833
$ws.Range("L80").Value = 4235.9375
$ws.Range("M80").Value = -1110.5833
$ws.Range("N80").Value = -6231.9375
$ws.Range("H81").Value = 495011800
$ws.Range("J81").Value = 495011800
$ws.Range("L81").Value = 495011800
$ws.Range("N81").Value = -495013796
$ws.Range("H83").Value = 3324.2144
$ws.Range("I83").Value = 2108.5833
$ws.Range("J83").Value = 4235.9375
$ws.Range("K83").Value = 10542.9165
$ws.Range("L83").Value = 21179.6875
$ws.Range("M83").Value = -5550.916499999999
$ws.Range("N83").Value = -31163.6875
$ws.Range("H84").Value = 495011800
$ws.Range("J84").Value = 495011800
$ws.Range("L84").Value = 1485035400
$ws.Range("N84").Value = -1485045384
$ws.Range("H109").Value = 40000
$ws.Range("I109").Value = 40000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 40000
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("M109").Value = -38960

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1938.4736
$ws.Range("I55").Value = 1325.909
$ws.Range("K55").Value = 1325.909
$ws.Range("M55").Value = -1152.909
$ws.Range("H93").Value = 4575.8184
$ws.Range("I93").Value = 1933.3334
$ws.Range("K93").Value = 1933.3334
$ws.Range("M93").Value = -685.3334
$ws.Range("H122").Value = 3491.1035
$ws.Range("I122").Value = 3259.8076
$ws.Range("J122").Value = 5495.6665
$ws.Range("K122").Value = 9779.4228
$ws.Range("L122").Value = 16486.9995
$ws.Range("M122").Value = -7329.4228
$ws.Range("N122").Value = -21386.9995
$ws.Range("H132").Value = 4278.2144
$ws.Range("I132").Value = 3537
$ws.Range("K132").Value = 10611
$ws.Range("M132").Value = -8081

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 12380
$ws.Range("J33").Value = 12380
$ws.Range("L33").Value = 12380
$ws.Range("N33").Value = -12880
$ws.Range("H36").Value = 12380
$ws.Range("J36").Value = 12380
$ws.Range("L36").Value = 12380
$ws.Range("N36").Value = -12880
$ws.Range("H40").Value = 28333.334
$ws.Range("J40").Value = 28333.334
$ws.Range("L40").Value = 28333.334
$ws.Range("N40").Value = -28631.334
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H122").Value = 3940.4827
$ws.Range("I122").Value = 2098.913
$ws.Range("K122").Value = 6296.739
$ws.Range("M122").Value = -3846.739
$ws.Range("H132").Value = 12347.143
$ws.Range("I132").Value = 13297.105
$ws.Range("K132").Value = 39891.315
$ws.Range("M132").Value = -37361.315
